$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.570.12'
$ws.Range("E2").Value = '  +3.90%  '

# Row 3
$ws.Range("D3").Value = '3.510.84'
$ws.Range("E3").Value = '  +3.86%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
${tmpStyle} = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.28'
$ws.Range("D5").Style = ${tmpStyle}
$ws.Range("E5").Value = '  +6.11%  '

# Row 6
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
${tmpStyle} = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.17'
$ws.Range("D6").Style = ${tmpStyle}
$ws.Range("E6").Value = '  +6.67%  '

# Row 7
${tmpStyle} = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").Style = ${tmpStyle}
$ws.Range("E7").Value = '  +5.21%  '

# Row 8
$ws.Range("D8").Value = '3.507.78'
$ws.Range("E8").Value = '  +3.80%  '

# Row 9
$ws.Range("E9").Value = '  +0.00%  '

# Row 10
${tmpStyle} = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.631'
$ws.Range("D10").Style = ${tmpStyle}
$ws.Range("E10").Value = '  +3.36%  '

# Row 11
${tmpStyle} = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.154'
$ws.Range("D11").Style = ${tmpStyle}
$ws.Range("E11").Value = '  +14.27%  '

# Row 12
${tmpStyle} = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.49'
$ws.Range("D12").Style = ${tmpStyle}
$ws.Range("E12").Value = '  +1.58%  '

# Row 13
${tmpStyle} = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000272'
$ws.Range("D13").Style = ${tmpStyle}
$ws.Range("E13").Value = '  +5.56%  '

# Row 14
${tmpStyle} = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.30'
$ws.Range("D14").Style = ${tmpStyle}
$ws.Range("E14").Value = '  +1.99%  '

# Row 15
$ws.Range("D15").Value = '4.081.73'
$ws.Range("E15").Value = '  +3.98%  '

# Row 16
$ws.Range("D16").Value = '3.519.73'
$ws.Range("E16").Value = '  +3.96%  '

# Row 17
${tmpStyle} = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.61'
$ws.Range("D17").Style = ${tmpStyle}
$ws.Range("E17").Value = '  +5.65%  '

# Row 18
${tmpStyle} = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.121'
$ws.Range("D18").Style = ${tmpStyle}
$ws.Range("E18").Value = '  +2.90%  '

# Row 19
$ws.Range("D19").Value = '66.686.51'
$ws.Range("E19").Value = '  +4.23%  '

# Row 20
${tmpStyle} = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.03'
$ws.Range("D20").Style = ${tmpStyle}
$ws.Range("E20").Value = '  +6.77%  '

# Row 21
${tmpStyle} = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.994'
$ws.Range("D21").Style = ${tmpStyle}
$ws.Range("E21").Value = '  +2.94%  '

# Row 22
${tmpStyle} = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '421.50'
$ws.Range("D22").Style = ${tmpStyle}
$ws.Range("E22").Value = '  +12.14%  '

# Row 23
${tmpStyle} = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.07'
$ws.Range("D23").Style = ${tmpStyle}
$ws.Range("E23").Value = '  +10.16%  '

# Row 24
${tmpStyle} = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.00'
$ws.Range("D24").Style = ${tmpStyle}
$ws.Range("E24").Value = '  +5.35%  '

# Row 25
${tmpStyle} = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.15'
$ws.Range("D25").Style = ${tmpStyle}
$ws.Range("E25").Value = '  -2.14%  '

# Row 26
${tmpStyle} = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.98'
$ws.Range("D26").Style = ${tmpStyle}
$ws.Range("E26").Value = '  -5.13%  '

# Row 27
${tmpStyle} = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.90'
$ws.Range("D27").Style = ${tmpStyle}
$ws.Range("E27").Value = '  +6.95%  '

# Row 28
${tmpStyle} = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.09'
$ws.Range("D28").Style = ${tmpStyle}
$ws.Range("E28").Value = '  -1.01%  '

# Row 29
${tmpStyle} = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '12.24'
$ws.Range("D29").Style = ${tmpStyle}
$ws.Range("E29").Value = '  +7.99%  '

# Row 30
${tmpStyle} = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.07'
$ws.Range("D30").Style = ${tmpStyle}
$ws.Range("E30").Value = '  +10.07%  '

# Row 31
${tmpStyle} = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.19'
$ws.Range("D31").Style = ${tmpStyle}
$ws.Range("E31").Value = '  +4.03%  '

# Row 32
${tmpStyle} = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '629.32'
$ws.Range("D32").Style = ${tmpStyle}
$ws.Range("E32").Value = '  -0.24%  '

# Row 33
${tmpStyle} = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.58'
$ws.Range("D33").Style = ${tmpStyle}
$ws.Range("E33").Value = '  +1.84%  '

# Row 34
${tmpStyle} = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.71'
$ws.Range("D34").Style = ${tmpStyle}
$ws.Range("E34").Value = '  +4.23%  '

# Row 35
$ws.Range("E35").Value = '  +4.44%  '

# Row 36
${tmpStyle} = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '60.19'
$ws.Range("D36").Style = ${tmpStyle}
$ws.Range("E36").Value = '  +3.65%  '

# Row 37
$ws.Range("D37").Value = '0.0₃0826'
$ws.Range("E37").Value = '  +11.62%  '

# Row 38
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
${tmpStyle} = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.148'
$ws.Range("D38").Style = ${tmpStyle}
$ws.Range("E38").Value = '  +18.23%  '

# Row 39
$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
${tmpStyle} = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.95'
$ws.Range("D39").Style = ${tmpStyle}
$ws.Range("E39").Value = '  +4.41%  '

# Row 40
$ws.Range("E40").Value = '  -0.06%  '

# Row 41
${tmpStyle} = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.385'
$ws.Range("D41").Style = ${tmpStyle}
$ws.Range("E41").Value = '  +1.17%  '

# Row 42
$ws.Range("E42").Value = '  +11.96%  '

# Row 43
$ws.Range("D43").Value = '3.121.18'
$ws.Range("E43").Value = '  +4.83%  '

# Row 44
${tmpStyle} = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = ${tmpStyle}
$ws.Range("E44").Value = '  +0.09%  '

# Row 45
${tmpStyle} = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.62'
$ws.Range("D45").Style = ${tmpStyle}
$ws.Range("E45").Value = '  -2.94%  '

# Row 46
${tmpStyle} = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.85'
$ws.Range("D46").Style = ${tmpStyle}
$ws.Range("E46").Value = '  +9.05%  '

# Row 47
${tmpStyle} = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.32'
$ws.Range("D47").Style = ${tmpStyle}
$ws.Range("E47").Value = '  +9.64%  '

# Row 48
${tmpStyle} = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0415'
$ws.Range("D48").Style = ${tmpStyle}
$ws.Range("E48").Value = '  +4.67%  '

# Row 49
$ws.Range("E49").Value = '  +2.62%  '

# Row 50
${tmpStyle} = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.132'
$ws.Range("D50").Style = ${tmpStyle}
$ws.Range("E50").Value = '  +5.56%  '

# Row 51
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
${tmpStyle} = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '139.06'
$ws.Range("D51").Style = ${tmpStyle}
$ws.Range("E51").Value = '  +2.16%  '
